$d = $word.ActiveDocument

# 1) First occurrence: "... CPF nº {{ num_cpf }}, RG nº {{ num_rg }}, declaro ..."
#    "RG nº" becomes "documento de identidade {{sigla_identidade}} "
$d.Content.Find.Execute(
    "RG nº",            # FindText
    $true,               # MatchCase
    $false,               # MatchWholeWord
    $false,               # MatchWildcards
    $false,               # MatchSoundsLike
    $false,               # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "documento de identidade {{sigla_identidade}} ",  # ReplaceWith
    2                     # Replace (wdReplaceAll)
) | Out-Null

# 2) Second occurrence: "RG: {{ num_rg }}" becomes "{{ sigla_identidade }}: {{ num_rg }}"
$d.Content.Find.Execute(
    "RG:",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "{{ sigla_identidade }}:",
    2
) | Out-Null
